$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at E:F (existing E..H shift right to G..J),
# carrying formatting from the neighboring column (D) the same way Excel does.
$ws.Range("E:F").Insert()

# New "moneda" column (E) - header keeps the inherited header style,
# the value cell needs the same style D2 uses (the "data" row style).
$ws.Range("E1").Value = "moneda"
$ws.Range("E2").Value = "USD"
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New "cuenta" column (F) - this column does not carry over a copied style.
$ws.Range("F1").Value = "cuenta"
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value = "Cuenta Simple"

# The new "moneda" column is narrower (fits "USD"/"moneda"); "cuenta" reuses
# the width the old "Estado" column used to have.
$ws.Columns("E:E").ColumnWidth = 6.8
$ws.Columns("F:F").ColumnWidth = 11.8

# Update the shifted transaction data with the new run's values
$ws.Range("H2").Value = "AAACT231810953MV37 "
$ws.Range("I2").Value = "30 jun. 2023, 14:07:30"

# J2 used to hold a plain number; the new run logs the account number as text
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "1010826248"
$ws.Range("J2").Style = "Normal"

# Update the selection to match the saved view
$ws.Range("G8").Select()
